# Auto-generated Word COM-interop script to apply the RFP draft edits.
$d = $word.ActiveDocument

# Paragraph 3
$p = $d.Paragraphs(3)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Our company''s full legal name is XYZ Capital Management, Inc.'

# Paragraph 5
$p = $d.Paragraphs(5)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Yes, our company has a comprehensive Code of Conduct and Ethics that guides our operations and business practices. This code is rigorously enforced to ensure integrity, transparency, and accountability in all our activities. It outlines clear standards for professional behavior and ethical decision-making, which are mandatory for all employees and executives. Additionally, we conduct regular training sessions and audits to ensure compliance and continuous awareness of ethical standards. This framework not only supports our commitment to lawful and ethical conduct but also enhances our corporate governance and fosters a culture of trust and respect among our stakeholders.'

# Paragraph 7
$p = $d.Paragraphs(7)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Our company maintains a strict policy prohibiting facilitation payments to public officials or any other parties. This policy aligns with our commitment to ethical business practices and compliance with anti-corruption laws globally. We enforce this stance through rigorous training programs for all employees and a robust compliance framework that includes monitoring and enforcement mechanisms to ensure adherence to these standards. Additionally, our internal policies are regularly reviewed and updated to respond to new legal and regulatory developments, ensuring our practices remain transparent and uphold the highest levels of integrity.'

# Paragraph 9
$p = $d.Paragraphs(9)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Our proposed technology solution is designed to optimize operational efficiency and enhance decision-making through advanced analytics and machine learning. Leveraging state-of-the-art software and hardware, our platform integrates seamlessly with existing systems, ensuring a smooth transition and minimal disruption to ongoing operations.' + `
    [char]11 + `
    [char]11 + `
    'The core of our solution is a robust data analytics engine that utilizes real-time processing to deliver actionable insights. This engine supports a wide range of functionalities including predictive analytics, risk assessment, and portfolio management, tailored to meet the specific needs of your organization.' + `
    [char]11 + `
    [char]11 + `
    'Security is paramount in our technology architecture. We employ industry-leading encryption methods and comprehensive data protection protocols to safeguard sensitive information against unauthorized access and cyber threats.' + `
    [char]11 + `
    [char]11 + `
    'Additionally, our solution is scalable and flexible, designed to grow with your organization. It can easily adapt to new market conditions and evolving business requirements, ensuring long-term sustainability and cost-effectiveness.' + `
    [char]11 + `
    [char]11 + `
    'Our team of experts will provide ongoing support and training, ensuring that your staff are fully equipped to maximize the benefits of our technology. This approach not only enhances user adoption but also ensures that you achieve the highest return on your investment.'

# Paragraph 11
$p = $d.Paragraphs(11)
$r = $p.Range
$r.End = $r.End - 1
$r.Text = 'Our proposed technology solution offers a suite of advantages designed to enhance operational efficiencies and provide superior performance compared to in-house built systems. Key benefits include:' + `
    [char]11 + `
    [char]11 + `
    '1. **Scalability and Flexibility**: Our solution is built to scale seamlessly with your business growth without the need for significant additional investments in infrastructure. This flexibility allows you to adapt more dynamically to market changes and business needs compared to more rigid, in-house systems.' + `
    [char]11 + `
    [char]11 + `
    '2. **Cost Efficiency**: By leveraging our technology, your firm can significantly reduce both initial and ongoing expenses associated with developing, maintaining, and updating an in-house system. Our solution eliminates the need for costly hardware acquisitions and reduces the burden on your IT staff, allowing them to focus on more strategic initiatives.' + `
    [char]11 + `
    [char]11 + `
    '3. **Advanced Security Features**: We employ state-of-the-art security protocols that are continually updated to respond to emerging threats, providing a level of security typically beyond the reach of in-house solutions. This includes comprehensive data encryption, regular security audits, and compliance with international standards.' + `
    [char]11 + `
    [char]11 + `
    '4. **Enhanced Data Analytics and Reporting**: Our system integrates advanced analytics tools that provide actionable insights and detailed reporting capabilities. This enables better decision-making based on real-time data, a capability that in-house systems often struggle to match due to limitations in processing power and analytical features.' + `
    [char]11 + `
    [char]11 + `
    '5. **Continuous Updates and Innovation**: Unlike in-house systems that may suffer from slower update cycles, our solution receives continuous updates that incorporate the latest technological advancements and regulatory changes. This ensures that your operations always stay ahead of industry curves and compliance requirements.' + `
    [char]11 + `
    [char]11 + `
    '6. **Expert Support and Maintenance**: Our team of experts is available around the clock to provide support and perform system maintenance, ensuring high availability and minimal downtime. This level of professional support is often challenging and costly to replicate with in-house resources.' + `
    [char]11 + `
    [char]11 + `
    '7. **Proven Track Record**: Our technology has been successfully implemented across multiple firms, demonstrating tangible improvements in operational efficiency, cost reduction, and risk management. Clients have reported enhanced processing speeds, reduced error rates, and improved client satisfaction.' + `
    [char]11 + `
    [char]11 + `
    'In summary, our technology solution not only addresses the typical challenges faced by in-house systems but also provides a robust framework for scalability, security, and continuous improvement, leading to measurable enhancements in operational efficiency and strategic advantage.'
